# Aggiornamento fino a 21 marzo: aggiunge le righe mancanti con i nuovi dati giornalieri.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date seriali (colonna A) per i nuovi giorni da aggiungere, a partire dalla riga 230.
$dates = @(44304, 44305, 44306, 44307)
$startRow = 230

$row = $startRow
foreach ($d in $dates) {
    $ws.Range("A$row").Value = $d
    $ws.Range("B$row").Value = 0
    $ws.Range("C$row").Value = 0
    $ws.Range("D$row").Value = 0
    $row = $row + 1
}
$endRow = $row - 1

# Applica alle nuove celle lo stesso formato delle celle dell'ultima riga esistente (229),
# cosi' la colonna A mantiene lo stile/formato data e le colonne B:D restano col formato standard.
$ws.Range("A229").Copy()
$ws.Range("A${startRow}:A${endRow}").PasteSpecial(-4122)

$ws.Range("B229:D229").Copy()
$ws.Range("B${startRow}:D${endRow}").PasteSpecial(-4122)

$excel.CutCopyMode = 0
